$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 697; everything from old row 697 onward
# shifts down by two (old 697->699, ..., old 731->733).
$ws.Rows("697:698").Insert()

# New row 697 (Primera, fecha 2022-07-11)
$ws.Range("A697").Value = 3
$ws.Range("B697").Value = "Femacal de La Calera"
$ws.Range("C697").Value = "Coquimbo"
$ws.Range("D697").Value = 44753
$ws.Range("E697").Value = 5
$ws.Range("F697").Value = 100112023
$ws.Range("G697").Value = "Brócoli"
$ws.Range("H697").Value = "Sin especificar"
$ws.Range("I697").Value = "Primera"
$ws.Range("J697").Value = 1900
$ws.Range("K697").Value = 800
$ws.Range("L697").Value = 900
$ws.Range("M697").Value = 850
$ws.Range("N697").Value = "$/unidad"
$ws.Range("O697").Value = "Provincia de Quillota"
$ws.Range("P697").Value = 850
$ws.Range("Q697").Value = 1
$ws.Range("R697").Value = "Hortaliza"

# New row 698 (Segunda, fecha 2022-07-11)
$ws.Range("A698").Value = 3
$ws.Range("B698").Value = "Femacal de La Calera"
$ws.Range("C698").Value = "Coquimbo"
$ws.Range("D698").Value = 44753
$ws.Range("E698").Value = 5
$ws.Range("F698").Value = 100112023
$ws.Range("G698").Value = "Brócoli"
$ws.Range("H698").Value = "Sin especificar"
$ws.Range("I698").Value = "Segunda"
$ws.Range("J698").Value = 900
$ws.Range("K698").Value = 700
$ws.Range("L698").Value = 700
$ws.Range("M698").Value = 700
$ws.Range("N698").Value = "$/unidad"
$ws.Range("O698").Value = "Provincia de Quillota"
$ws.Range("P698").Value = 700
$ws.Range("Q698").Value = 1
$ws.Range("R698").Value = "Hortaliza"
